$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (39) down through
# the new rows (40-49) so column A / E keep their styles (s="1" / s="2").
$ws.Range("A39:V39").Copy()
$ws.Range("A40:V49").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 40 (Indice 39)
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "kuwait"
$ws.Cells.Item(40, 3).Value = "premier-league"
$ws.Cells.Item(40, 4).Value = "2023-2024"
$ws.Cells.Item(40, 5).Value = 45263.69791666666
$ws.Cells.Item(40, 6).Value = "Al Arabi"
$ws.Cells.Item(40, 7).Value = 3
$ws.Cells.Item(40, 8).Value = "Al Kuwait"
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 3.55
$ws.Cells.Item(40, 11).Value = "03/12/2023 05:12"
$ws.Cells.Item(40, 12).Value = 3.34
$ws.Cells.Item(40, 13).Value = "03/12/2023 15:20"
$ws.Cells.Item(40, 14).Value = 3.59
$ws.Cells.Item(40, 15).Value = "03/12/2023 05:12"
$ws.Cells.Item(40, 16).Value = 3.46
$ws.Cells.Item(40, 17).Value = "03/12/2023 15:20"
$ws.Cells.Item(40, 18).Value = 1.88
$ws.Cells.Item(40, 19).Value = "03/12/2023 05:12"
$ws.Cells.Item(40, 20).Value = 2.01
$ws.Cells.Item(40, 21).Value = "03/12/2023 15:20"
$ws.Cells.Item(40, 22).Value = "https://www.betexplorer.com/football/kuwait/premier-league/al-arabi-kuwait-al-kuwait/4S9ME3fU/"

# Row 41 (Indice 40)
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "kuwait"
$ws.Cells.Item(41, 3).Value = "premier-league"
$ws.Cells.Item(41, 4).Value = "2023-2024"
$ws.Cells.Item(41, 5).Value = 45265.63888888889
$ws.Cells.Item(41, 6).Value = "Khaitan"
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = "Al Jahra"
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 2.58
$ws.Cells.Item(41, 11).Value = "05/12/2023 03:42"
$ws.Cells.Item(41, 12).Value = 2.26
$ws.Cells.Item(41, 13).Value = "05/12/2023 14:24"
$ws.Cells.Item(41, 14).Value = 3.26
$ws.Cells.Item(41, 15).Value = "05/12/2023 03:42"
$ws.Cells.Item(41, 16).Value = 3.37
$ws.Cells.Item(41, 17).Value = "05/12/2023 14:41"
$ws.Cells.Item(41, 18).Value = 2.51
$ws.Cells.Item(41, 19).Value = "05/12/2023 03:42"
$ws.Cells.Item(41, 20).Value = 2.88
$ws.Cells.Item(41, 21).Value = "05/12/2023 14:24"
$ws.Cells.Item(41, 22).Value = "https://www.betexplorer.com/football/kuwait/premier-league/khaitan-al-jahra/6eWiUrfH/"

# Row 42 (Indice 41)
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "kuwait"
$ws.Cells.Item(42, 3).Value = "premier-league"
$ws.Cells.Item(42, 4).Value = "2023-2024"
$ws.Cells.Item(42, 5).Value = 45265.75694444445
$ws.Cells.Item(42, 6).Value = "Al Shabab"
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = "Al Qadisiya"
$ws.Cells.Item(42, 9).Value = 4
$ws.Cells.Item(42, 10).Value = 4.16
$ws.Cells.Item(42, 11).Value = "05/12/2023 06:12"
$ws.Cells.Item(42, 12).Value = 4.45
$ws.Cells.Item(42, 13).Value = "05/12/2023 17:12"
$ws.Cells.Item(42, 14).Value = 3.55
$ws.Cells.Item(42, 15).Value = "05/12/2023 06:12"
$ws.Cells.Item(42, 16).Value = 3.66
$ws.Cells.Item(42, 17).Value = "05/12/2023 17:12"
$ws.Cells.Item(42, 18).Value = 1.75
$ws.Cells.Item(42, 19).Value = "05/12/2023 06:12"
$ws.Cells.Item(42, 20).Value = 1.7
$ws.Cells.Item(42, 21).Value = "05/12/2023 17:12"
$ws.Cells.Item(42, 22).Value = "https://www.betexplorer.com/football/kuwait/premier-league/al-shabab-al-qadisiya/CYWmVOvB/"

# Row 43 (Indice 42)
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "kuwait"
$ws.Cells.Item(43, 3).Value = "premier-league"
$ws.Cells.Item(43, 4).Value = "2023-2024"
$ws.Cells.Item(43, 5).Value = 45266.69791666666
$ws.Cells.Item(43, 6).Value = "Al Salmiya"
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = "Al-Fahaheel"
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 1.83
$ws.Cells.Item(43, 11).Value = "06/12/2023 05:11"
$ws.Cells.Item(43, 12).Value = 1.51
$ws.Cells.Item(43, 13).Value = "06/12/2023 16:44"
$ws.Cells.Item(43, 14).Value = 3.69
$ws.Cells.Item(43, 15).Value = "06/12/2023 05:11"
$ws.Cells.Item(43, 16).Value = 4.27
$ws.Cells.Item(43, 17).Value = "06/12/2023 16:44"
$ws.Cells.Item(43, 18).Value = 3.63
$ws.Cells.Item(43, 19).Value = "06/12/2023 05:11"
$ws.Cells.Item(43, 20).Value = 5.31
$ws.Cells.Item(43, 21).Value = "06/12/2023 16:36"
$ws.Cells.Item(43, 22).Value = "https://www.betexplorer.com/football/kuwait/premier-league/al-salmiya-al-fahaheel/YsYqW4P4/"

# Row 44 (Indice 43)
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "kuwait"
$ws.Cells.Item(44, 3).Value = "premier-league"
$ws.Cells.Item(44, 4).Value = "2023-2024"
$ws.Cells.Item(44, 5).Value = 45267.63888888889
$ws.Cells.Item(44, 6).Value = "Al Kuwait"
$ws.Cells.Item(44, 7).Value = 2
$ws.Cells.Item(44, 8).Value = "Al Naser"
$ws.Cells.Item(44, 9).Value = 1
$ws.Cells.Item(44, 10).Value = 1.43
$ws.Cells.Item(44, 11).Value = "07/12/2023 03:42"
$ws.Cells.Item(44, 12).Value = 1.52
$ws.Cells.Item(44, 13).Value = "07/12/2023 14:50"
$ws.Cells.Item(44, 14).Value = 4.56
$ws.Cells.Item(44, 15).Value = "07/12/2023 03:42"
$ws.Cells.Item(44, 16).Value = 4.35
$ws.Cells.Item(44, 17).Value = "07/12/2023 14:50"
$ws.Cells.Item(44, 18).Value = 5.69
$ws.Cells.Item(44, 19).Value = "07/12/2023 03:42"
$ws.Cells.Item(44, 20).Value = 5.03
$ws.Cells.Item(44, 21).Value = "07/12/2023 14:50"
$ws.Cells.Item(44, 22).Value = "https://www.betexplorer.com/football/kuwait/premier-league/al-kuwait-al-naser/4tkeT29N/"

# Row 45 (Indice 44)
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "kuwait"
$ws.Cells.Item(45, 3).Value = "premier-league"
$ws.Cells.Item(45, 4).Value = "2023-2024"
$ws.Cells.Item(45, 5).Value = 45267.75694444445
$ws.Cells.Item(45, 6).Value = "Al Arabi"
$ws.Cells.Item(45, 7).Value = 4
$ws.Cells.Item(45, 8).Value = "Kazma SC"
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 1.95
$ws.Cells.Item(45, 11).Value = "07/12/2023 06:12"
$ws.Cells.Item(45, 12).Value = 1.98
$ws.Cells.Item(45, 13).Value = "07/12/2023 18:09"
$ws.Cells.Item(45, 14).Value = 3.62
$ws.Cells.Item(45, 15).Value = "07/12/2023 06:12"
$ws.Cells.Item(45, 16).Value = 3.55
$ws.Cells.Item(45, 17).Value = "07/12/2023 18:09"
$ws.Cells.Item(45, 18).Value = 3.3
$ws.Cells.Item(45, 19).Value = "07/12/2023 06:12"
$ws.Cells.Item(45, 20).Value = 3.35
$ws.Cells.Item(45, 21).Value = "07/12/2023 18:09"
$ws.Cells.Item(45, 22).Value = "https://www.betexplorer.com/football/kuwait/premier-league/al-arabi-kuwait-kazma-sc/n5v0SMOT/"

# Row 46 (Indice 45)
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "kuwait"
$ws.Cells.Item(46, 3).Value = "premier-league"
$ws.Cells.Item(46, 4).Value = "2023-2024"
$ws.Cells.Item(46, 5).Value = 45274.64236111111
$ws.Cells.Item(46, 6).Value = "Khaitan"
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = "Al Qadisiya"
$ws.Cells.Item(46, 9).Value = 3
$ws.Cells.Item(46, 10).Value = 8.02
$ws.Cells.Item(46, 11).Value = "14/12/2023 03:42"
$ws.Cells.Item(46, 12).Value = 7.66
$ws.Cells.Item(46, 13).Value = "14/12/2023 14:51"
$ws.Cells.Item(46, 14).Value = 4.99
$ws.Cells.Item(46, 15).Value = "14/12/2023 03:42"
$ws.Cells.Item(46, 16).Value = 4.9
$ws.Cells.Item(46, 17).Value = "14/12/2023 14:17"
$ws.Cells.Item(46, 18).Value = 1.31
$ws.Cells.Item(46, 19).Value = "14/12/2023 03:42"
$ws.Cells.Item(46, 20).Value = 1.33
$ws.Cells.Item(46, 21).Value = "14/12/2023 13:45"
$ws.Cells.Item(46, 22).Value = "https://www.betexplorer.com/football/kuwait/premier-league/khaitan-al-qadisiya/hjPrE8Is/"

# Row 47 (Indice 46)
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "kuwait"
$ws.Cells.Item(47, 3).Value = "premier-league"
$ws.Cells.Item(47, 4).Value = "2023-2024"
$ws.Cells.Item(47, 5).Value = 45274.75
$ws.Cells.Item(47, 6).Value = "Al Shabab"
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = "Al-Fahaheel"
$ws.Cells.Item(47, 9).Value = 2
$ws.Cells.Item(47, 10).Value = 2.19
$ws.Cells.Item(47, 11).Value = "14/12/2023 06:12"
$ws.Cells.Item(47, 12).Value = 2.03
$ws.Cells.Item(47, 13).Value = "14/12/2023 16:30"
$ws.Cells.Item(47, 14).Value = 3.39
$ws.Cells.Item(47, 15).Value = "14/12/2023 06:12"
$ws.Cells.Item(47, 16).Value = 3.42
$ws.Cells.Item(47, 17).Value = "14/12/2023 16:30"
$ws.Cells.Item(47, 18).Value = 2.94
$ws.Cells.Item(47, 19).Value = "14/12/2023 06:12"
$ws.Cells.Item(47, 20).Value = 3.33
$ws.Cells.Item(47, 21).Value = "14/12/2023 16:04"
$ws.Cells.Item(47, 22).Value = "https://www.betexplorer.com/football/kuwait/premier-league/al-shabab-al-fahaheel/nDOnDSXm/"

# Row 48 (Indice 47)
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "kuwait"
$ws.Cells.Item(48, 3).Value = "premier-league"
$ws.Cells.Item(48, 4).Value = "2023-2024"
$ws.Cells.Item(48, 5).Value = 45275.58333333334
$ws.Cells.Item(48, 6).Value = "Al Jahra"
$ws.Cells.Item(48, 7).Value = 2
$ws.Cells.Item(48, 8).Value = "Al Kuwait"
$ws.Cells.Item(48, 9).Value = 3
$ws.Cells.Item(48, 10).Value = 8.82
$ws.Cells.Item(48, 11).Value = "15/12/2023 02:12"
$ws.Cells.Item(48, 12).Value = 7.98
$ws.Cells.Item(48, 13).Value = "15/12/2023 13:59"
$ws.Cells.Item(48, 14).Value = 5.34
$ws.Cells.Item(48, 15).Value = "15/12/2023 02:12"
$ws.Cells.Item(48, 16).Value = 5.22
$ws.Cells.Item(48, 17).Value = "15/12/2023 13:59"
$ws.Cells.Item(48, 18).Value = 1.26
$ws.Cells.Item(48, 19).Value = "15/12/2023 02:12"
$ws.Cells.Item(48, 20).Value = 1.3
$ws.Cells.Item(48, 21).Value = "15/12/2023 13:58"
$ws.Cells.Item(48, 22).Value = "https://www.betexplorer.com/football/kuwait/premier-league/al-jahra-al-kuwait/0GSjCnmf/"

# Row 49 (Indice 48)
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "kuwait"
$ws.Cells.Item(49, 3).Value = "premier-league"
$ws.Cells.Item(49, 4).Value = "2023-2024"
$ws.Cells.Item(49, 5).Value = 45275.69791666666
$ws.Cells.Item(49, 6).Value = "Kazma SC"
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = "Al Naser"
$ws.Cells.Item(49, 9).Value = 2
$ws.Cells.Item(49, 10).Value = 2.2
$ws.Cells.Item(49, 11).Value = "15/12/2023 05:13"
$ws.Cells.Item(49, 12).Value = 2.21
$ws.Cells.Item(49, 13).Value = "15/12/2023 16:17"
$ws.Cells.Item(49, 14).Value = 3.44
$ws.Cells.Item(49, 15).Value = "15/12/2023 05:13"
$ws.Cells.Item(49, 16).Value = 3.58
$ws.Cells.Item(49, 17).Value = "15/12/2023 16:17"
$ws.Cells.Item(49, 18).Value = 2.88
$ws.Cells.Item(49, 19).Value = "15/12/2023 05:13"
$ws.Cells.Item(49, 20).Value = 2.82
$ws.Cells.Item(49, 21).Value = "15/12/2023 16:17"
$ws.Cells.Item(49, 22).Value = "https://www.betexplorer.com/football/kuwait/premier-league/kazma-sc-al-naser/ruIeB620/"
